$d = $word.ActiveDocument

# The ACS Number Request form previously pointed readers to
# "acstnrequest@microsoft.com". The correct TNS alias is
# "acstns@microsoft.com" - update the visible hyperlink text while
# leaving the mailto: target relationship untouched.
$updated = $false
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.TextToDisplay -eq "acstnrequest@microsoft.com") {
        $h.TextToDisplay = "acstns@microsoft.com"
        $updated = $true
    }
}

if (-not $updated) {
    # Fallback: directly edit the text range if the hyperlink lookup above
    # did not find a match for some reason.
    $rng = $d.Range(0, $d.Content.End)
    if ($rng.Find.Execute("acstnrequest@microsoft.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $start = $rng.Start
        $middle = $d.Range($start + 5, $start + 12)
        $middle.Text = "s"
    }
}

Write-Output "updated: $updated"
